$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-11 (columns D, M, N, O, P, R, S) per the commit diff.
# This is a permutation/reshuffle of the weekly price-report rows.
$data = @{
    2  = @{ D = 44585; M = 160; N = 6500; O = 7000; P = 6750; R = "Provincia de Curicó"; S = 3375 }
    3  = @{ D = 44209; M = 58;  N = 6000; O = 6000; P = 6000; R = "Provincia de Curicó"; S = 3000 }
    4  = @{ D = 44211; M = 45;  N = 6000; O = 6000; P = 6000; R = "Provincia de Curicó"; S = 3000 }
    5  = @{ D = 44587; M = 165; N = 6500; O = 7000; P = 6742; R = "Provincia de Linares"; S = 3371 }
    6  = @{ D = 44586; M = 80;  N = 7000; O = 7000; P = 7000; R = "Provincia de Curicó"; S = 3500 }
    7  = @{ D = 44589; M = 60;  N = 6000; O = 6000; P = 6000; R = "Provincia de Curicó"; S = 3000 }
    8  = @{ D = 44592; M = 30;  N = 8000; O = 8000; P = 8000; R = "Provincia de Linares"; S = 4000 }
    9  = @{ D = 44582; M = 150; N = 6000; O = 6500; P = 6233; R = "Provincia de Curicó"; S = 3116 }
    10 = @{ D = 44588; M = 160; N = 6500; O = 7000; P = 6750; R = "Provincia de Curicó"; S = 3375 }
    11 = @{ D = 44214; M = 48;  N = 6000; O = 6000; P = 6000; R = "Provincia de Linares"; S = 3000 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("R$row").Value = $vals.R
    $ws.Range("S$row").Value = $vals.S
}
